$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new row 36 first so it inherits the formatting (style) of row 35
# directly above it before any new content is placed in that row.
$ws.Rows("36:36").Insert(-4121, -4163)
$ws.Range("A36").Value = "producto comodin"
$ws.Range("A37").Value = "cambiar precio de producto"

# Row 6: add Responsable "Lucas" and percent 80%
$ws.Range("B6").Value = "Lucas"
$ws.Range("C6").Value = 0.8
$ws.Range("C6").NumberFormat = "0%"

# Row 20: add percent 100%
$ws.Range("C20").Value = 1
$ws.Range("C20").NumberFormat = "0%"

# Row 25: add Responsable "Agustina"
$ws.Range("B25").Value = "Agustina"

# Row 35: add Responsable "Agustina " (trailing space) and percent 100%
$ws.Range("B35").Value = "Agustina "
$ws.Range("C35").Value = 1
$ws.Range("C35").NumberFormat = "0%"

# Update the active selection to reflect where the user ended up editing
$ws.Range("C36").Select()
